# Refined metadata to be additional tab
#
# 1) Refresh the "time_taken" column (F2:F23) on the existing "data" sheet
#    with the new run's timestamps.
# 2) Add a new "metadata" worksheet (after "data") describing the panel
#    query that produced the data sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1) Update time_taken values on the "data" sheet -----------------------
$timeTaken = @(
    "2021-10-05 14:35:08.188498",
    "2021-10-05 14:35:08.188506",
    "2021-10-05 14:35:08.188509",
    "2021-10-05 14:35:08.188512",
    "2021-10-05 14:35:08.188515",
    "2021-10-05 14:35:08.188518",
    "2021-10-05 14:35:08.188520",
    "2021-10-05 14:35:08.188523",
    "2021-10-05 14:35:08.188526",
    "2021-10-05 14:35:08.188528",
    "2021-10-05 14:35:08.188531",
    "2021-10-05 14:35:08.188534",
    "2021-10-05 14:35:08.188536",
    "2021-10-05 14:35:08.188539",
    "2021-10-05 14:35:08.188541",
    "2021-10-05 14:35:08.188544",
    "2021-10-05 14:35:08.188547",
    "2021-10-05 14:35:08.188550",
    "2021-10-05 14:35:08.188552",
    "2021-10-05 14:35:08.188555",
    "2021-10-05 14:35:08.188557",
    "2021-10-05 14:35:08.188560"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws1.Range("F$row").Value = $timeTaken[$i]
}

# --- 2) Add the "metadata" sheet, placed right after "data" ----------------
$meta = $wb.Worksheets.Add($null, $ws1)
$meta.Name = "metadata"

# Header row (B1:G1) — bold, centered, thin-bordered, matching the "data"
# sheet's header style.
$headerRange = $meta.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Row 2 — the index cell (A2) carries the same header-like style as the
# "data" sheet's index column.
$a2 = $meta.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$meta.Range("B2").Value = "Osteopetrosis"
$meta.Range("C2").Value = 150

# data_version ("0.8") must stay a text value, not be coerced to a number.
$d2 = $meta.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "0.8"
$d2.Style = $ws1.Range("B2").Style

$meta.Range("E2").Value = "2021-05-14T06:29:27.208328Z"
$meta.Range("F2").Value = "2021-10-05 14:35:08.184864"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/150/?format=json"

# Restore "data" as the active sheet/selection (the diff only adds the new
# tab; it does not change which sheet is active).
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
